# "Refined metadata to be additional tab"
# 1. Refresh the time_taken values on the existing "data" sheet.
# 2. Add a new "metadata" worksheet right after "data" and populate it with
#    the panel-query metadata (name/id/version/timestamps/request url).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh the time_taken values on the data sheet --------------------
$timeTaken = @(
    "2021-10-05 14:22:09.029035",
    "2021-10-05 14:22:09.029043",
    "2021-10-05 14:22:09.029046",
    "2021-10-05 14:22:09.029049",
    "2021-10-05 14:22:09.029051",
    "2021-10-05 14:22:09.029054",
    "2021-10-05 14:22:09.029057",
    "2021-10-05 14:22:09.029059",
    "2021-10-05 14:22:09.029062",
    "2021-10-05 14:22:09.029065",
    "2021-10-05 14:22:09.029067",
    "2021-10-05 14:22:09.029070",
    "2021-10-05 14:22:09.029072",
    "2021-10-05 14:22:09.029075",
    "2021-10-05 14:22:09.029078",
    "2021-10-05 14:22:09.029080",
    "2021-10-05 14:22:09.029083",
    "2021-10-05 14:22:09.029085"
)

for ($i = 0; $i -lt $timeTaken.Count; $i++) {
    $row = $i + 2
    $dataSheet.Range("F$row").Value = $timeTaken[$i]
}

# --- 2. Add the "metadata" worksheet right after "data" --------------------
# Duplicate "data" (so sheetPr/sheetFormatPr/pageMargins/outline settings and
# the existing header/border styling all carry over for free), rename it,
# then wipe its contents and fill in the metadata table.
$dataSheet.Copy($null, $dataSheet)
$metaSheet = $wb.Worksheets.Item("data (2)")
$metaSheet.Name = "metadata"

$metaSheet.Cells.Clear()

# Re-apply the bold/bordered/centered header style (copied from the "data"
# sheet's own header row) to the new header row, and the bordered index-cell
# style to A2.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$metaSheet.Application.CutCopyMode = $false

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Pancreatitis"
$metaSheet.Range("C2").Value = 386
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "2.10"
$metaSheet.Range("D2").Style = "Normal"
$metaSheet.Range("E2").Value = "2021-04-12T10:13:42.285010Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:22:09.025886"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/386/?format=json"

# Keep "data" as the active/selected tab (the diff only appends the new
# sheet to the workbook's sheet list; it doesn't change which tab is active).
$dataSheet.Activate()

Write-Output "done"
